$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'69.784.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.10%  "

# Row 3
$ws.Range("D3").Value = "'3.500.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.91%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'598.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.50%  "

# Row 6
$ws.Range("D6").Value = "'194.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.29%  "

# Row 7
$ws.Range("D7").Value = "'0.622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.58%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").Value = "'0.207"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.39%  "

# Row 10
$ws.Range("D10").Value = "'0.650"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.54%  "

# Row 11
$ws.Range("D11").Value = "'53.73"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.38%  "

# Row 12
$ws.Range("D12").Value = "'0.0000300"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.34%  "

# Row 13
$ws.Range("D13").Value = "'9.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.77%  "

# Row 14
$ws.Range("D14").Value = "'4.056.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.84%  "

# Row 15
$ws.Range("D15").Value = "'606.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.62%  "

# Row 16
$ws.Range("D16").Value = "'69.939.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.10%  "

# Row 17
$ws.Range("D17").Value = "'18.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.75%  "

# Row 18
$ws.Range("D18").Value = "'12.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.39%  "

# Row 19
$ws.Range("D19").Value = "'3.500.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.86%  "

# Row 20
$ws.Range("D20").Value = "'0.120"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.60%  "

# Row 21
$ws.Range("D21").Value = "'0.990"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.13%  "

# Row 22
$ws.Range("D22").Value = "'18.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.49%  "

# Row 23
$ws.Range("D23").Value = "'104.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +10.91%  "

# Row 24
$ws.Range("D24").Value = "'5.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.68%  "

# Row 25
$ws.Range("E25").Value = "  -2.14%  "

# Row 26
$ws.Range("E26").Value = "  +3.02%  "

# Row 27
$ws.Range("D27").Value = "'10.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "

# Row 28
$ws.Range("D28").Value = "'9.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.89%  "

# Row 29
$ws.Range("D29").Value = "'33.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.42%  "

# Row 30
$ws.Range("D30").Value = "'4.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +25.24%  "

# Row 31
$ws.Range("D31").Value = "'7.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.59%  "

# Row 32
$ws.Range("D32").Value = "'12.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.07%  "

# Row 33
$ws.Range("E33").Value = "  +1.42%  "

# Row 34
$ws.Range("D34").Value = "'63.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.34%  "

# Row 35
$ws.Range("D35").Value = "'3.725.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.11%  "

# Row 36
$ws.Range("D36").Value = "'0.0₃0808"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.96%  "

# Row 37
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "

# Row 38
$ws.Range("E38").Value = "  -7.22%  "

# Row 39
$ws.Range("E39").Value = "  -2.66%  "

# Row 40
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'36.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.91%  "

# Row 41
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'3.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.06%  "

# Row 42
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "'500.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.00%  "

# Row 43
$ws.Range("E43").Value = "  +0.44%  "

# Row 44
$ws.Range("D44").Value = "'0.0457"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.33%  "

# Row 45
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'3.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.54%  "

# Row 46
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.140"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.11%  "

# Row 47
$ws.Range("E47").Value = "  -3.43%  "

# Row 48
$ws.Range("D48").Value = "'1.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.37%  "

# Row 49
$ws.Range("E49").Value = "  -3.51%  "

# Row 50
$ws.Range("D50").Value = "'131.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.70%  "

# Row 51
$ws.Range("D51").Value = "'0.000241"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.23%  "
